$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.42681884765625
$ws.Range("B1").Value = 1.379914283752441
$ws.Range("C1").Value = 3.436492443084717
$ws.Range("D1").Value = 2.743883848190308
$ws.Range("E1").Value = 0.874331533908844
